$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (single decimal point, e.g. "245.82") are forced to Text format first,
# then the style is reset back to Normal so no lingering number format remains.

$ws.Range("D2").Value = '35.838.37'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").Value = '1.989.96'
$ws.Range("E3").Value = '  -3.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.641'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.53'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.53%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.33'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.366'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.957'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '2.278.18'
$ws.Range("E15").Value = '  -3.54%  '
$ws.Range("E16").Value = '  -3.15%  '
$ws.Range("D17").Value = '1.992.04'
$ws.Range("E17").Value = '  -3.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.77%  '
$ws.Range("D19").Value = '35.763.08'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '0.0₃0852'
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.11%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.74%  '
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.86%  '
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0985'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.67%  '
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.77%  '
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.36%  '
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0957'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("D48").Value = '1.373.60'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.80%  '
